$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Values: B1 and A2 are numeric 0, B2 is the text label.
$ws.Range("B1").Value = 0
$ws.Range("A2").Value = 0
$ws.Range("B2").Value = "disconnected_elements"

# Build the full style (bold font, thin box border, centered/top aligned)
# on B1 first ...
$rng1 = $ws.Range("B1")
$rng1.Borders.LineStyle = 1
$rng1.Borders.Weight = 2
$rng1.Font.Bold = $true
$rng1.HorizontalAlignment = -4108
$rng1.VerticalAlignment = -4160

# ... then clone that exact formatting onto A2 via copy/paste-special so
# both cells end up sharing the very same cell style (xf) entry instead of
# each accumulating their own separate style records.
$rng1.Copy()
$ws.Range("A2").PasteSpecial(-4122)
$excel.CutCopyMode = $false
